$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.971.62'
$ws.Range("E2").Value = '  +2.05%  '
$ws.Range("D3").Value = '1.655.69'
$ws.Range("E3").Value = '  +2.71%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'215.14"
$ws.Range("E5").Value = '  +1.50%  '
$ws.Range("D6").Value = "'0.508"
$ws.Range("E6").Value = '  +2.33%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +2.66%  '
$ws.Range("D10").Value = "'20.12"
$ws.Range("E10").Value = '  +4.61%  '
$ws.Range("E11").Value = '  +3.93%  '
$ws.Range("D12").Value = '1.889.74'
$ws.Range("E12").Value = '  +2.72%  '
$ws.Range("D13").Value = '1.655.79'
$ws.Range("E13").Value = '  +2.65%  '
$ws.Range("E14").Value = '  +2.21%  '
$ws.Range("E15").Value = '  +2.94%  '
$ws.Range("D16").Value = "'65.31"
$ws.Range("E16").Value = '  +2.75%  '
$ws.Range("D17").Value = '26.978.85'
$ws.Range("E17").Value = '  +2.06%  '
$ws.Range("D18").Value = "'236.90"
$ws.Range("E18").Value = '  +0.86%  '
$ws.Range("D19").Value = '0.0₃0739'
$ws.Range("E19").Value = '  +1.90%  '
$ws.Range("D20").Value = "'7.78"
$ws.Range("E20").Value = '  +1.52%  '
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("D22").Value = "'4.43"
$ws.Range("E22").Value = '  +3.76%  '
$ws.Range("D23").Value = "'9.30"
$ws.Range("E23").Value = '  +2.87%  '
$ws.Range("D24").Value = "'2.22"
$ws.Range("E24").Value = '  +0.94%  '
$ws.Range("D25").Value = "'145.29"
$ws.Range("E25").Value = '  -0.91%  '
$ws.Range("E26").Value = '  +2.12%  '
$ws.Range("E27").Value = '  +0.83%  '
$ws.Range("D28").Value = "'15.87"
$ws.Range("E28").Value = '  +2.71%  '
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("E30").Value = '  +0.42%  '
$ws.Range("E31").Value = '  +1.82%  '
$ws.Range("D32").Value = '1.554.85'
$ws.Range("E32").Value = '  +3.56%  '
$ws.Range("E33").Value = '  +2.03%  '
$ws.Range("D34").Value = "'3.08"
$ws.Range("E34").Value = '  +4.47%  '
$ws.Range("E35").Value = '  +7.97%  '
$ws.Range("E36").Value = '  -0.28%  '
$ws.Range("D37").Value = "'0.582"
$ws.Range("E37").Value = '  +3.75%  '
$ws.Range("E38").Value = '  +9.72%  '
$ws.Range("E39").Value = '  +2.73%  '
$ws.Range("D40").Value = "'6.04"
$ws.Range("E40").Value = '  +3.96%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D42").Value = "'66.38"
$ws.Range("E42").Value = '  +7.68%  '
$ws.Range("D43").Value = "'0.976"
$ws.Range("E43").Value = '  +5.48%  '
$ws.Range("E44").Value = '  +2.50%  '
$ws.Range("D45").Value = '1.798.04'
$ws.Range("E45").Value = '  +2.63%  '
$ws.Range("E46").Value = '  +1.77%  '
$ws.Range("D47").Value = "'90.18"
$ws.Range("E47").Value = '  +0.58%  '
$ws.Range("E48").Value = '  +2.99%  '
$ws.Range("D49").Value = "'0.0997"
$ws.Range("E49").Value = '  +3.89%  '
$ws.Range("E50").Value = '  +1.00%  '
$ws.Range("D51").Value = "'7.68"
$ws.Range("E51").Value = '  +2.86%  '
